$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-07"

# Update the label for the September row
$ws.Range("A10").Value = "September (through 09-07)"

# Update September row (row 10) values for each year column (B-I)
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 17
$ws.Range("G10").Value = 24
$ws.Range("H10").Value = 29
$ws.Range("I10").Value = 35

# Update Total row (row 11) values for each year column (B-I)
$ws.Range("B11").Value = 201
$ws.Range("C11").Value = 394
$ws.Range("D11").Value = 568
$ws.Range("E11").Value = 499
$ws.Range("F11").Value = 366
$ws.Range("G11").Value = 808
$ws.Range("H11").Value = 1099
$ws.Range("I11").Value = 1172
